$d = $word.ActiveDocument

# 1) Title: "1 What is Blender" -> "The Shear Tool"
$p1 = $d.Paragraphs(1).Range
$p1.Find.Execute("1 What is Blender", $true, $false, $false, $false, $false, $true, 1, $false, "The Shear Tool", 2)

# 2) Keywords: insert "The Shear Tool, " before "Blender, 3D Modeling, Animation, Graphic Art"
$p6 = $d.Paragraphs(6).Range
$p6.Find.Execute("Blender, 3D Modeling, Animation, Graphic Art", $true, $false, $false, $false, $false, $true, 1, $false, "The Shear Tool, Blender, 3D Modeling, Animation, Graphic Art", 2)

# 3) Description: replace the "what the 3D modeling program ... about." sentence
$p9 = $d.Paragraphs(9).Range
$p9.Find.Execute("what the 3D modeling program “Blender " + [char]34 + " is all about.", $true, $false, $false, $false, $false, $true, 1, $false, "how to use the Shear tool, in Edit Mode, while working in the Blender application", 2)

# 4) Category: insert "The Shear Tool, " before "Blender, 3D Modeling, Animation, Graphic Art"
$p11 = $d.Paragraphs(11).Range
$p11.Find.Execute("Blender, 3D Modeling, Animation, Graphic Art", $true, $false, $false, $false, $false, $true, 1, $false, "The Shear Tool, Blender, 3D Modeling, Animation, Graphic Art", 2)

# 5) Revised date
$p18 = $d.Paragraphs(18).Range
$p18.Find.Execute("Wednesday, December 11, 2024", $true, $false, $false, $false, $false, $true, 1, $false, "Monday, January 27, 2025", 2)

# 6) URL
$p20 = $d.Paragraphs(20).Range
$p20.Find.Execute("Enlightenment/Articles/2024/8-Blender-2024/1-What-Is-Blender/1-What-Is-Blender.html", $true, $false, $false, $false, $false, $true, 1, $false, "Enlightenment/Articles/2025/1-Blender-Continued/2-Edit-Mode/1-The-Menus/1-The-Tools-Menu/12-The-Shear-Tool/The-Shear-Tool.html", 2)
